# Updates the cryptos list values (price + 1h volume change) on Sheet1, rows 2-51.
# Each entry is (row, col, newValue, isNumericLooking). Generated from the
# authoritative cell-level diff of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    ,@(2, 4, '43.209.12', $false)
    ,@(2, 5, '  +2.75%  ', $false)
    ,@(3, 4, '2.310.05', $false)
    ,@(3, 5, '  +4.48%  ', $false)
    ,@(4, 4, '1.00', $true)
    ,@(4, 5, '  -0.04%  ', $false)
    ,@(5, 4, '253.08', $true)
    ,@(5, 5, '  +0.60%  ', $false)
    ,@(6, 4, '0.644', $true)
    ,@(6, 5, '  +3.60%  ', $false)
    ,@(7, 4, '74.67', $true)
    ,@(7, 5, '  +10.38%  ', $false)
    ,@(8, 4, '1.00', $true)
    ,@(8, 5, '  -0.17%  ', $false)
    ,@(9, 4, '0.654', $true)
    ,@(9, 5, '  +6.41%  ', $false)
    ,@(10, 4, '39.62', $true)
    ,@(10, 5, '  +2.32%  ', $false)
    ,@(11, 4, '0.0993', $true)
    ,@(11, 5, '  +6.11%  ', $false)
    ,@(12, 4, '59.39', $true)
    ,@(12, 5, '  +0.05%  ', $false)
    ,@(13, 4, '7.41', $true)
    ,@(13, 5, '  +5.66%  ', $false)
    ,@(14, 5, '  +1.50%  ', $false)
    ,@(15, 4, '2.649.29', $false)
    ,@(15, 5, '  +4.18%  ', $false)
    ,@(16, 4, '15.55', $true)
    ,@(16, 5, '  +7.98%  ', $false)
    ,@(17, 4, '0.882', $true)
    ,@(17, 5, '  +1.63%  ', $false)
    ,@(18, 4, '2.308.84', $false)
    ,@(18, 5, '  +4.35%  ', $false)
    ,@(19, 4, '43.078.88', $false)
    ,@(19, 5, '  +2.73%  ', $false)
    ,@(20, 5, '  +5.14%  ', $false)
    ,@(21, 5, '  +3.37%  ', $false)
    ,@(22, 4, '72.87', $true)
    ,@(22, 5, '  +0.83%  ', $false)
    ,@(23, 4, '234.00', $true)
    ,@(23, 5, '  +1.24%  ', $false)
    ,@(24, 5, '  +10.30%  ', $false)
    ,@(25, 4, '3.92', $true)
    ,@(25, 5, '  +1.10%  ', $false)
    ,@(26, 4, '11.67', $true)
    ,@(26, 5, '  +4.62%  ', $false)
    ,@(27, 5, '  -0.05%  ', $false)
    ,@(28, 5, '  +1.34%  ', $false)
    ,@(29, 4, '3.64', $true)
    ,@(29, 5, '  -1.29%  ', $false)
    ,@(30, 5, '  -0.19%  ', $false)
    ,@(31, 4, '167.66', $true)
    ,@(31, 5, '  +0.49%  ', $false)
    ,@(32, 4, '21.21', $true)
    ,@(32, 5, '  +4.02%  ', $false)
    ,@(33, 4, '6.41', $true)
    ,@(33, 5, '  +8.81%  ', $false)
    ,@(34, 5, '  +5.75%  ', $false)
    ,@(35, 4, '0.0820', $true)
    ,@(35, 5, '  +5.37%  ', $false)
    ,@(36, 4, '32.26', $true)
    ,@(36, 5, '  +24.14%  ', $false)
    ,@(37, 5, '  +4.02%  ', $false)
    ,@(38, 4, '4.66', $true)
    ,@(38, 5, '  +13.69%  ', $false)
    ,@(39, 5, '  +4.61%  ', $false)
    ,@(40, 5, '  -0.86%  ', $false)
    ,@(41, 4, '14.50', $true)
    ,@(41, 5, '  +21.54%  ', $false)
    ,@(42, 4, '2.36', $true)
    ,@(42, 5, '  +6.14%  ', $false)
    ,@(43, 4, '6.02', $true)
    ,@(43, 5, '  +6.78%  ', $false)
    ,@(44, 5, '  +10.42%  ', $false)
    ,@(45, 2, 'FraxShare', $false)
    ,@(45, 3, 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', $false)
    ,@(45, 4, '9.18', $true)
    ,@(45, 5, '  +7.53%  ', $false)
    ,@(46, 2, 'MultiversX', $false)
    ,@(46, 3, 'https://coinranking.com/coin/omwkOTglq+multiversx-egld', $false)
    ,@(46, 4, '62.39', $true)
    ,@(46, 5, '  +1.91%  ', $false)
    ,@(47, 4, '4.91', $true)
    ,@(47, 5, '  -4.74%  ', $false)
    ,@(48, 5, '  +4.18%  ', $false)
    ,@(49, 2, 'ARBITRUM', $false)
    ,@(49, 3, 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', $false)
    ,@(49, 4, '1.19', $true)
    ,@(49, 5, '  +3.75%  ', $false)
    ,@(50, 2, 'BinanceUSD', $false)
    ,@(50, 3, 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', $false)
    ,@(50, 4, '1.00', $true)
    ,@(50, 5, '  +0.09%  ', $false)
    ,@(51, 4, '98.94', $true)
    ,@(51, 5, '  +6.68%  ', $false)
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $val = $u[2]
    $isNumericLooking = $u[3]
    $cell = $ws.Cells.Item($row, $col)

    if ($isNumericLooking) {
        # This text value looks like a plain number (e.g. "1.00", "234.00",
        # "0.0993"). Left alone, Excel auto-converts such input to the Number
        # type, silently dropping significant trailing zeros. The source data
        # stores it as text, so force the Text format before assigning, then
        # restore the default "Normal" style so no stray formatting diff is
        # introduced.
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}
